$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Octubre de 2020 a las 10:04"

# Row 57 - Chequia: updated stats
$ws.Range("B57").Value = 70834
$ws.Range("C57").Value = 63
$ws.Range("D57").Value = 33557
$ws.Range("E57").Value = 36619
$ws.Range("G57").Value = 3
$ws.Range("H57").Value = 658

# Row 59 - Singapur: updated stats
$ws.Range("B59").Value = 57786
$ws.Range("C59").Value = 21
$ws.Range("E59").Value = 271

# Row 64 - Armenia: updated stats
$ws.Range("B64").Value = 50850
$ws.Range("C64").Value = 491
$ws.Range("D64").Value = 44219
$ws.Range("E64").Value = 5668
$ws.Range("G64").Value = 4
$ws.Range("H64").Value = 963

# Rows 80-81: Hungria moves above Australia, Hungria gets new stats,
# Australia keeps the former Hungria row's old stats (values shift down a row)
$ws.Range("A80").Value = "Hungria"
$ws.Range("B80").Value = 27309
$ws.Range("C80").Value = 848
$ws.Range("D80").Value = 6118
$ws.Range("E80").Value = 20410
$ws.Range("G80").Value = 16
$ws.Range("H80").Value = 781

$ws.Range("A81").Value = "Australia"
$ws.Range("B81").Value = 27096
$ws.Range("C81").Value = 18
$ws.Range("D81").Value = 24784
$ws.Range("E81").Value = 1424
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 888

# Rows 115-116: Georgia moves above Jamaica, Georgia gets new stats,
# Jamaica keeps the former Georgia row's old stats
$ws.Range("A115").Value = "Georgia"
$ws.Range("B115").Value = 6640
$ws.Range("C115").Value = 448
$ws.Range("D115").Value = 3419
$ws.Range("E115").Value = 3182
$ws.Range("H115").Value = 39

$ws.Range("A116").Value = "Jamaica"
$ws.Range("B116").Value = 6482
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 1867
$ws.Range("E116").Value = 4508
$ws.Range("H116").Value = 107

# Rows 142-143: Estonia moves above Sri Lanka, Estonia gets new stats,
# Sri Lanka keeps the former Estonia row's old stats
$ws.Range("A142").Value = "Estonia"
$ws.Range("B142").Value = 3450
$ws.Range("C142").Value = 81
$ws.Range("D142").Value = 2643
$ws.Range("E142").Value = 743
$ws.Range("H142").Value = 64

$ws.Range("A143").Value = "Sri Lanka"
$ws.Range("B143").Value = 3380
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 3230
$ws.Range("E143").Value = 137
$ws.Range("H143").Value = 13

$wb.Save()
